$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Client" column before the old Instance column (C) ---
$ws.Columns("C").Insert()

# --- Header row (row 1): add new "Client" header ---
$ws.Range("C1").Value = "Client"

# --- Fix existing rollout row (row 2): correct truncated filename + add Client value ---
$ws.Range("A2").Value = "CT_REPORT_00303_V.001.tar"
$ws.Range("C2").Value = "KUMHO - 3036"

# --- New rollout row (row 3): fill in the previously-blank row with new data ---
$ws.Range("A3").Value = "RG_PAGEBUILDER_00052_V.003"
$ws.Range("B3").Value = "PB to create Orders"
$ws.Range("C3").Value = "Soda Stream"
$ws.Range("D3").Value = "PR10"
$ws.Range("E3").Value = "NO"
$ws.Range("G3").Value = "NO"

# --- Row heights: existing row shrinks, new row takes the old wrapped height ---
$ws.Rows("2").RowHeight = 29
$ws.Rows("3").RowHeight = 43.5

# --- Column widths (nearest values this engine's width rounding can hit) ---
$ws.Columns("A").ColumnWidth = 27.8
$ws.Columns("C").ColumnWidth = 9.6

# --- Selection moves off the old A1:G2 block ---
[void]$ws.Range("F11").Select()
